$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 723.8182
$ws.Cells.Item(28, 9).Value = 663.3333
$ws.Cells.Item(28, 10).Value = 796.4
$ws.Cells.Item(28, 11).Value = 663.3333
$ws.Cells.Item(28, 12).Value = 796.4
$ws.Cells.Item(28, 13).Value = -178.3333
$ws.Cells.Item(28, 14).Value = -1766.4

$ws.Cells.Item(125, 8).Value = 905.4
$ws.Cells.Item(125, 9).Value = 894.8889
$ws.Cells.Item(125, 11).Value = 8054.0001
$ws.Cells.Item(125, 13).Value = -5594.0001

$ws.Cells.Item(131, 8).Value = 24466.455
$ws.Cells.Item(131, 9).Value = 31208.182
$ws.Cells.Item(131, 10).Value = 4241.273
$ws.Cells.Item(131, 11).Value = 93624.546
$ws.Cells.Item(131, 12).Value = 12723.819
$ws.Cells.Item(131, 13).Value = -88584.546
$ws.Cells.Item(131, 14).Value = -22803.819

$ws.Cells.Item(132, 8).Value = 1889.2727
$ws.Cells.Item(132, 9).Value = 1919.5555
$ws.Cells.Item(132, 10).Value = 254
$ws.Cells.Item(132, 11).Value = 5758.666499999999
$ws.Cells.Item(132, 12).Value = 762
$ws.Cells.Item(132, 13).Value = -3228.666499999999
$ws.Cells.Item(132, 14).Value = -5822

$ws.Cells.Item(135, 8).Value = 1168
$ws.Cells.Item(135, 9).Value = 901.8570999999999
$ws.Cells.Item(135, 10).Value = 1789
$ws.Cells.Item(135, 11).Value = 8116.7139
$ws.Cells.Item(135, 12).Value = 16101
$ws.Cells.Item(135, 13).Value = -5581.7139
$ws.Cells.Item(135, 14).Value = -21171

$ws.Cells.Item(137, 8).Value = 40486.348
$ws.Cells.Item(137, 9).Value = 1688.0769
$ws.Cells.Item(137, 10).Value = 79284.62
$ws.Cells.Item(137, 11).Value = 5064.2307
$ws.Cells.Item(137, 12).Value = 237853.86
$ws.Cells.Item(137, 13).Value = -2514.2307
$ws.Cells.Item(137, 14).Value = -242953.86

$ws.Cells.Item(138, 8).Value = 4099.515
$ws.Cells.Item(138, 9).Value = 1011
$ws.Cells.Item(138, 10).Value = 5087.84
$ws.Cells.Item(138, 11).Value = 3033
$ws.Cells.Item(138, 12).Value = 15263.52
$ws.Cells.Item(138, 13).Value = 2107
$ws.Cells.Item(138, 14).Value = -25543.52

$ws.Cells.Item(141, 8).Value = 996.8
$ws.Cells.Item(141, 9).Value = 996
$ws.Cells.Item(141, 10).Value = 1000
$ws.Cells.Item(141, 11).Value = 2988
$ws.Cells.Item(141, 12).Value = 3000
$ws.Cells.Item(141, 13).Value = 2192
$ws.Cells.Item(141, 14).Value = -13360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5695.778
$ws.Cells.Item(61, 9).Value = 2169.6
$ws.Cells.Item(61, 11).Value = 2169.6
$ws.Cells.Item(61, 13).Value = -1957.6

$ws.Cells.Item(74, 8).Value = 18989.5
$ws.Cells.Item(74, 9).Value = 21787.4
$ws.Cells.Item(74, 10).Value = 5000
$ws.Cells.Item(74, 11).Value = 21787.4
$ws.Cells.Item(74, 12).Value = 5000
$ws.Cells.Item(74, 13).Value = -20913.4
$ws.Cells.Item(74, 14).Value = -6748

$ws.Cells.Item(77, 8).Value = 18989.5
$ws.Cells.Item(77, 9).Value = 21787.4
$ws.Cells.Item(77, 10).Value = 5000
$ws.Cells.Item(77, 11).Value = 108937
$ws.Cells.Item(77, 12).Value = 25000
$ws.Cells.Item(77, 13).Value = -104569
$ws.Cells.Item(77, 14).Value = -33736

$ws.Cells.Item(132, 8).Value = 2193
$ws.Cells.Item(132, 9).Value = 1669.9546
$ws.Cells.Item(132, 10).Value = 3471.5557
$ws.Cells.Item(132, 11).Value = 5009.8638
$ws.Cells.Item(132, 12).Value = 10414.6671
$ws.Cells.Item(132, 13).Value = -2479.8638
$ws.Cells.Item(132, 14).Value = -15474.6671

$ws.Cells.Item(136, 8).Value = 5695.778
$ws.Cells.Item(136, 9).Value = 2169.6
$ws.Cells.Item(136, 11).Value = 6508.799999999999
$ws.Cells.Item(136, 13).Value = -3958.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1658.7778
$ws.Cells.Item(134, 9).Value = 1628.625
$ws.Cells.Item(134, 10).Value = 1900
$ws.Cells.Item(134, 11).Value = 4885.875
$ws.Cells.Item(134, 12).Value = 5700
$ws.Cells.Item(134, 13).Value = -2350.875
$ws.Cells.Item(134, 14).Value = -10770

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1009.63635
$ws.Cells.Item(22, 10).Value = 1300.6666
$ws.Cells.Item(22, 12).Value = 1300.6666
$ws.Cells.Item(22, 14).Value = -2000.6666

$ws.Cells.Item(31, 8).Value = 4927.5386
$ws.Cells.Item(31, 9).Value = 4039.3333
$ws.Cells.Item(31, 10).Value = 5688.857
$ws.Cells.Item(31, 11).Value = 4039.3333
$ws.Cells.Item(31, 12).Value = 5688.857
$ws.Cells.Item(31, 13).Value = -3744.3333
$ws.Cells.Item(31, 14).Value = -6278.857

$ws.Cells.Item(34, 8).Value = 4927.5386
$ws.Cells.Item(34, 9).Value = 4039.3333
$ws.Cells.Item(34, 10).Value = 5688.857
$ws.Cells.Item(34, 11).Value = 4039.3333
$ws.Cells.Item(34, 12).Value = 5688.857
$ws.Cells.Item(34, 13).Value = -3837.3333
$ws.Cells.Item(34, 14).Value = -6092.857

$ws.Cells.Item(107, 8).Value = 608.8889
$ws.Cells.Item(107, 9).Value = 525.55554
$ws.Cells.Item(107, 10).Value = 775.55554
$ws.Cells.Item(107, 11).Value = 525.55554
$ws.Cells.Item(107, 12).Value = 775.55554
$ws.Cells.Item(107, 13).Value = 1394.44446
$ws.Cells.Item(107, 14).Value = -4615.55554

$ws.Cells.Item(132, 8).Value = 2105.7827
$ws.Cells.Item(132, 9).Value = 1263.0555
$ws.Cells.Item(132, 10).Value = 5139.6
$ws.Cells.Item(132, 11).Value = 3789.1665
$ws.Cells.Item(132, 12).Value = 15418.8
$ws.Cells.Item(132, 13).Value = -1259.1665
$ws.Cells.Item(132, 14).Value = -20478.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 222.47058
$ws.Cells.Item(23, 9).Value = 197
$ws.Cells.Item(23, 11).Value = 591
$ws.Cells.Item(23, 13).Value = -356

$ws.Cells.Item(122, 8).Value = 8601.074000000001
$ws.Cells.Item(122, 9).Value = 11585.789
$ws.Cells.Item(122, 11).Value = 104272.101
$ws.Cells.Item(122, 13).Value = -101822.101

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 182.52942
$ws.Cells.Item(107, 9).Value = 200.21428
$ws.Cells.Item(107, 11).Value = 200.21428
$ws.Cells.Item(107, 13).Value = 1719.78572

$ws.Cells.Item(132, 8).Value = 3989.76
$ws.Cells.Item(132, 9).Value = 3848.2778
$ws.Cells.Item(132, 10).Value = 4353.5713
$ws.Cells.Item(132, 11).Value = 11544.8334
$ws.Cells.Item(132, 12).Value = 13060.7139
$ws.Cells.Item(132, 13).Value = -9014.8334
$ws.Cells.Item(132, 14).Value = -18120.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4089.182
$ws.Cells.Item(40, 9).Value = 3559
$ws.Cells.Item(40, 11).Value = 3559
$ws.Cells.Item(40, 13).Value = -3423

$ws.Cells.Item(122, 8).Value = 2071.5715
$ws.Cells.Item(122, 9).Value = 2105.15
$ws.Cells.Item(122, 10).Value = 1400
$ws.Cells.Item(122, 11).Value = 6315.450000000001
$ws.Cells.Item(122, 12).Value = 4200
$ws.Cells.Item(122, 13).Value = -3865.450000000001
$ws.Cells.Item(122, 14).Value = -9100

$ws.Cells.Item(136, 8).Value = 3269.2307
$ws.Cells.Item(136, 9).Value = 1856.4286
$ws.Cells.Item(136, 11).Value = 5569.2858
$ws.Cells.Item(136, 13).Value = -3019.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2598.5
$ws.Cells.Item(132, 9).Value = 2189.9333
$ws.Cells.Item(132, 10).Value = 3155.6365
$ws.Cells.Item(132, 11).Value = 6569.7999
$ws.Cells.Item(132, 12).Value = 9466.9095
$ws.Cells.Item(132, 13).Value = -4039.7999
$ws.Cells.Item(132, 14).Value = -14526.9095

$ws.Cells.Item(136, 8).Value = 14384.333
$ws.Cells.Item(136, 9).Value = 17637
$ws.Cells.Item(136, 11).Value = 52911
$ws.Cells.Item(136, 13).Value = -50361
